$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 878154
$ws.Range("I9").Value = 1148283.5
$ws.Range("K9").Value = 1148283.5
$ws.Range("M9").Value = -1148114.5
$ws.Range("H15").Value = 3073.923
$ws.Range("I15").Value = 3073.923
$ws.Range("K15").Value = 9221.769
$ws.Range("M15").Value = -9052.769
$ws.Range("H57").Value = 41629.715
$ws.Range("J57").Value = 41629.715
$ws.Range("L57").Value = 124889.145
$ws.Range("N57").Value = -125887.145
$ws.Range("H70").Value = 4775
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 4775
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 14325
$ws.Range("N70").Value = -14865
$ws.Range("H73").Value = 4775
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 4775
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 14325
$ws.Range("N73").Value = -16197
$ws.Range("H129").Value = 1144.4546
$ws.Range("J129").Value = 2612.7144
$ws.Range("L129").Value = 7838.1432
$ws.Range("N129").Value = -17838.1432
$ws.Range("H132").Value = 4125.3945
$ws.Range("I132").Value = 4278.222
$ws.Range("K132").Value = 12834.666
$ws.Range("M132").Value = -10304.666
$ws.Range("H137").Value = 1782.2
$ws.Range("I137").Value = 1808.3889
$ws.Range("J137").Value = 1546.5
$ws.Range("K137").Value = 5425.1667
$ws.Range("L137").Value = 4639.5
$ws.Range("M137").Value = -2875.1667
$ws.Range("N137").Value = -9739.5
$ws.Range("H138").Value = 4256.4116
$ws.Range("I138").Value = 1259.8334
$ws.Range("J138").Value = 5890.909
$ws.Range("K138").Value = 3779.5002
$ws.Range("L138").Value = 17672.727
$ws.Range("M138").Value = 1360.4998
$ws.Range("N138").Value = -27952.727
$ws.Range("N9").ClearContents()
$ws.Range("M70").ClearContents()
$ws.Range("M73").ClearContents()
$ws.Range("N132").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2290347.8
$ws.Range("I32").Value = 1160001
$ws.Range("K32").Value = 1160001
$ws.Range("M32").Value = -1159714
$ws.Range("H61").Value = 2662.5715
$ws.Range("I61").Value = 2280
$ws.Range("J61").Value = 3427.7144
$ws.Range("K61").Value = 2280
$ws.Range("L61").Value = 3427.7144
$ws.Range("M61").Value = -2068
$ws.Range("N61").Value = -3851.7144
$ws.Range("H102").Value = 2599.4285
$ws.Range("I102").Value = 1765.8889
$ws.Range("J102").Value = 4099.8
$ws.Range("K102").Value = 1765.8889
$ws.Range("L102").Value = 4099.8
$ws.Range("M102").Value = -143.8888999999999
$ws.Range("N102").Value = -7343.8
$ws.Range("H110").Value = 856.3333
$ws.Range("I110").Value = 700
$ws.Range("J110").Value = 934.5
$ws.Range("K110").Value = 700
$ws.Range("L110").Value = 934.5
$ws.Range("M110").Value = 1345
$ws.Range("N110").Value = -5024.5
$ws.Range("H112").Value = 29988.5
$ws.Range("J112").Value = 29988.5
$ws.Range("L112").Value = 29988.5
$ws.Range("N112").Value = -32942.5
$ws.Range("H132").Value = 1384.9474
$ws.Range("I132").Value = 1175.091
$ws.Range("J132").Value = 2770
$ws.Range("K132").Value = 3525.273
$ws.Range("L132").Value = 8310
$ws.Range("M132").Value = -995.2729999999997
$ws.Range("N132").Value = -13370
$ws.Range("H135").Value = 93737.60000000001
$ws.Range("J135").Value = 93737.60000000001
$ws.Range("L135").Value = 93737.60000000001
$ws.Range("N135").Value = -103877.6
$ws.Range("H136").Value = 2662.5715
$ws.Range("I136").Value = 2280
$ws.Range("J136").Value = 3427.7144
$ws.Range("K136").Value = 6840
$ws.Range("L136").Value = 10283.1432
$ws.Range("M136").Value = -4290
$ws.Range("N136").Value = -15383.1432
$ws.Range("N32").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2981.7568
$ws.Range("I86").Value = 3210.5454
$ws.Range("J86").Value = 2646.2
$ws.Range("K86").Value = 3210.5454
$ws.Range("L86").Value = 2646.2
$ws.Range("M86").Value = -2087.5454
$ws.Range("N86").Value = -4892.2
$ws.Range("H89").Value = 2981.7568
$ws.Range("I89").Value = 3210.5454
$ws.Range("J89").Value = 2646.2
$ws.Range("K89").Value = 16052.727
$ws.Range("L89").Value = 13231
$ws.Range("M89").Value = -10436.727
$ws.Range("N89").Value = -24463
$ws.Range("H94").Value = 222234320
$ws.Range("I94").Value = 333350500
$ws.Range("J94").Value = 1966.3334
$ws.Range("K94").Value = 333350500
$ws.Range("L94").Value = 1966.3334
$ws.Range("M94").Value = -333350049
$ws.Range("N94").Value = -2868.3334
$ws.Range("H134").Value = 1497.9429
$ws.Range("I134").Value = 901.11536
$ws.Range("J134").Value = 3222.111
$ws.Range("K134").Value = 2703.34608
$ws.Range("L134").Value = 9666.332999999999
$ws.Range("M134").Value = -168.3460800000003
$ws.Range("N134").Value = -14736.333
$ws.Range("H135").Value = 99985.5
$ws.Range("J135").Value = 99985.5
$ws.Range("L135").Value = 99985.5
$ws.Range("N135").Value = -110125.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 197.27272
$ws.Range("I7").Value = 195.375
$ws.Range("K7").Value = 195.375
$ws.Range("M7").Value = -82.375
$ws.Range("H15").Value = 1555
$ws.Range("J15").Value = 1555
$ws.Range("L15").Value = 1555
$ws.Range("N15").Value = -1895
$ws.Range("H22").Value = 406.66666
$ws.Range("I22").Value = 376.22223
$ws.Range("K22").Value = 376.22223
$ws.Range("M22").Value = -26.22223000000002
$ws.Range("H31").Value = 5684368.5
$ws.Range("I31").Value = 1546.4706
$ws.Range("K31").Value = 1546.4706
$ws.Range("M31").Value = -1251.4706
$ws.Range("H34").Value = 5684368.5
$ws.Range("I34").Value = 1546.4706
$ws.Range("K34").Value = 1546.4706
$ws.Range("M34").Value = -1344.4706
$ws.Range("H132").Value = 4211.107
$ws.Range("I132").Value = 3360.6
$ws.Range("K132").Value = 10081.8
$ws.Range("M132").Value = -7551.799999999999
$ws.Range("N7").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("N31").ClearContents()
$ws.Range("N34").ClearContents()
$ws.Range("N132").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 8398.5
$ws.Range("I56").Value = 8398.5
$ws.Range("K56").Value = 8398.5
$ws.Range("M56").Value = -7868.5
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 2453.4
$ws.Range("I2").Value = 2453.4
$ws.Range("K2").Value = 2453.4
$ws.Range("M2").Value = -2340.4
$ws.Range("H11").Value = 6265400
$ws.Range("J11").Value = 3248417
$ws.Range("L11").Value = 3248417
$ws.Range("N11").Value = -3248695
$ws.Range("H46").Value = 49973.5
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 49973.5
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 49973.5
$ws.Range("N46").Value = -50285.5
$ws.Range("H97").Value = 793
$ws.Range("I97").Value = 845.46155
$ws.Range("K97").Value = 845.46155
$ws.Range("M97").Value = -349.46155
$ws.Range("H132").Value = 1815.0416
$ws.Range("I132").Value = 1746.0646
$ws.Range("K132").Value = 5238.1938
$ws.Range("M132").Value = -2708.1938
$ws.Range("M46").ClearContents()
$ws.Range("N97").ClearContents()
$ws.Range("N132").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2374.6365
$ws.Range("I7").Value = 2210.375
$ws.Range("K7").Value = 2210.375
$ws.Range("M7").Value = -2098.375
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("H126").Value = 2374.6365
$ws.Range("I126").Value = 2210.375
$ws.Range("K126").Value = 6631.125
$ws.Range("M126").Value = -4161.125
$ws.Range("H132").Value = 4867
$ws.Range("I132").Value = 7066.3335
$ws.Range("K132").Value = 21199.0005
$ws.Range("M132").Value = -18669.0005
$ws.Range("N7").ClearContents()
$ws.Range("N110").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("N132").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 45027.5
$ws.Range("J38").Value = 49999
$ws.Range("L38").Value = 49999
$ws.Range("N38").Value = -50945
$ws.Range("H87").Value = 49999
$ws.Range("J87").Value = 49999
$ws.Range("L87").Value = 49999
$ws.Range("N87").Value = -52495
$ws.Range("H90").Value = 49999
$ws.Range("J90").Value = 49999
$ws.Range("L90").Value = 149997
$ws.Range("N90").Value = -162477
$ws.Range("H122").Value = 8334843
$ws.Range("I122").Value = 1367.92
$ws.Range("K122").Value = 4103.76
$ws.Range("M122").Value = -1653.76
$ws.Range("H132").Value = 4010.6875
$ws.Range("I132").Value = 3611.4666
$ws.Range("K132").Value = 10834.3998
$ws.Range("M132").Value = -8304.399800000001
$ws.Range("N122").ClearContents()
$ws.Range("N132").ClearContents()

Write-Host "Applied all Gilgamesh_Profits updates"